$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.805.37'
$ws.Range('E2').Value = '  -1.64%  '

# Row 3
$ws.Range('D3').Value = '1.892.00'
$ws.Range('E3').Value = '  -1.42%  '

# Row 4
$ws.Range('E4').Value = '  -0.24%  '

# Row 5
$ws.Range('D5').Value = '''0.7790'
$ws.Range('E5').Value = '  -3.56%  '

# Row 6
$ws.Range('D6').Value = '''244.12'
$ws.Range('E6').Value = '  -0.10%  '

# Row 7
$ws.Range('D7').Value = '''0.9999'
$ws.Range('E7').Value = '  -0.11%  '

# Row 8
$ws.Range('D8').Value = '''0.3126'
$ws.Range('E8').Value = '  -3.56%  '

# Row 9
$ws.Range('D9').Value = '''25.34'
$ws.Range('E9').Value = '  -6.78%  '

# Row 10
$ws.Range('D10').Value = '''0.07188'
$ws.Range('E10').Value = '  +1.33%  '

# Row 11
$ws.Range('D11').Value = '''0.08075'
$ws.Range('E11').Value = '  -0.27%  '

# Row 12
$ws.Range('D12').Value = '''0.7658'
$ws.Range('E12').Value = '  -2.26%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.932.06'
$ws.Range('E13').Value = '  +0.46%  '

# Row 14
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '''5.467'
$ws.Range('E14').Value = '  +1.04%  '

# Row 15
$ws.Range('D15').Value = '''92.40'
$ws.Range('E15').Value = '  -2.58%  '

# Row 16
$ws.Range('D16').Value = '''6.174'
$ws.Range('E16').Value = '  +2.45%  '

# Row 17
$ws.Range('D17').Value = '29.824.37'
$ws.Range('E17').Value = '  -1.59%  '

# Row 18
$ws.Range('D18').Value = '''13.95'
$ws.Range('E18').Value = '  -2.52%  '

# Row 19
$ws.Range('D19').Value = '''243.45'
$ws.Range('E19').Value = '  -2.25%  '

# Row 20
$ws.Range('D20').Value = '''0.000007770'
$ws.Range('E20').Value = '  -0.66%  '

# Row 21
$ws.Range('E21').Value = '  -0.12%  '

# Row 22
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').Value = '''8.147'
$ws.Range('E22').Value = '  +2.88%  '

# Row 23
$ws.Range('B23').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D23').Value = '2.140.92'
$ws.Range('E23').Value = '  -1.36%  '

# Row 24
$ws.Range('D24').Value = '''0.9999'
$ws.Range('E24').Value = '  -0.30%  '

# Row 25
$ws.Range('D25').Value = '''0.1575'
$ws.Range('E25').Value = '  -2.78%  '

# Row 26
$ws.Range('D26').Value = '''9.394'
$ws.Range('E26').Value = '  -0.80%  '

# Row 27
$ws.Range('D27').Value = '''162.14'
$ws.Range('E27').Value = '  -3.29%  '

# Row 28
$ws.Range('D28').Value = '''18.74'
$ws.Range('E28').Value = '  -1.94%  '

# Row 29
$ws.Range('D29').Value = '''2.050'
$ws.Range('E29').Value = '  -3.43%  '

# Row 30
$ws.Range('D30').Value = '''1.425'
$ws.Range('E30').Value = '  +3.78%  '

# Row 31
$ws.Range('D31').Value = '''1.549'
$ws.Range('E31').Value = '  +0.78%  '

# Row 32
$ws.Range('D32').Value = '''4.472'
$ws.Range('E32').Value = '  +2.67%  '

# Row 33
$ws.Range('D33').Value = '''4.102'
$ws.Range('E33').Value = '  -0.69%  '

# Row 34
$ws.Range('D34').Value = '''0.05527'
$ws.Range('E34').Value = '  -2.13%  '

# Row 35
$ws.Range('D35').Value = '''1.261'
$ws.Range('E35').Value = '  -3.14%  '

# Row 36
$ws.Range('D36').Value = '''0.7473'
$ws.Range('E36').Value = '  +0.79%  '

# Row 37
$ws.Range('D37').Value = '''1.006'
$ws.Range('E37').Value = '  +0.61%  '

# Row 38
$ws.Range('D38').Value = '''2.626'
$ws.Range('E38').Value = '  -3.34%  '

# Row 39
$ws.Range('D39').Value = '''0.01915'
$ws.Range('E39').Value = '  -1.74%  '

# Row 40
$ws.Range('D40').Value = '''2.775'
$ws.Range('E40').Value = '  -1.54%  '

# Row 41
$ws.Range('D41').Value = '1.136.84'
$ws.Range('E41').Value = '  +8.56%  '

# Row 42
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '''73.77'
$ws.Range('E42').Value = '  -0.06%  '

# Row 43
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = '''0.4428'
$ws.Range('E43').Value = '  -1.05%  '

# Row 44
$ws.Range('D44').Value = '''5.888'
$ws.Range('E44').Value = '  -1.52%  '

# Row 45
$ws.Range('D45').Value = '''0.8499'
$ws.Range('E45').Value = '  -0.48%  '

# Row 46
$ws.Range('D46').Value = '''103.98'
$ws.Range('E46').Value = '  +0.85%  '

# Row 47
$ws.Range('D47').Value = '''0.9995'
$ws.Range('E47').Value = '  -0.15%  '

# Row 48
$ws.Range('D48').Value = '''1.889'
$ws.Range('E48').Value = '  -2.12%  '

# Row 49
$ws.Range('D49').Value = '''9.988'
$ws.Range('E49').Value = '  +0.14%  '

# Row 50
$ws.Range('B50').Value = 'SynthetixNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D50').Value = '''3.039'
$ws.Range('E50').Value = '  +11.65%  '

# Row 51
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').Value = '''7.466'
$ws.Range('E51').Value = '  -2.17%  '
